$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("HomePage")
$ws2 = $wb.Worksheets.Item("LoginPage")

# Add the new object rows on the LoginPage sheet
$ws2.Range("A3").Value = "ddlSelectName"
$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = "'//select[@id='userSelect']"
$ws2.Range("D3").Value = "Dropdown"

$ws2.Range("A4").Value = "btnLogin"
$ws2.Range("B4").Value = 3
$ws2.Range("C4").Value = "'//button[@type='submit']"
$ws2.Range("D4").Value = "Button"

# Update selections/active sheet to match the final workbook view state
$ws1.Activate() | Out-Null
$ws1.Range("A3").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("D5").Select() | Out-Null
